$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Order matters for shared-string table construction - follow the same
# sequence the author used while updating progress notes.
$ws.Range("D3").Value = "No"
$ws.Range("E3").Value = "not available in apps-list for Agave apps"

$ws.Range("C2").Value = "dongwang/"

$ws.Range("B4").Value = "dave"
$ws.Range("D4").Value = "No"
$ws.Range("E4").Value = "version 2 is not public, working on getting published"

$ws.Range("D11").Value = "No"
$ws.Range("E11").Value = "Might be a problem with not being public, might be a problem with not sourcing dir correclty"
$ws.Range("C11").Value = "Plink"

# Column E width / best fit (author widened + auto-fit this column after
# lengthening its longest entry)
$ws.Columns.Item(5).ColumnWidth = 78.3

# Selection change
$ws.Range("C3").Select()
